# Applies the "ANOVA done and new diagrams" edit to the single-slide
# deck: the big pasted diagram group (Shape 1 on Slide 1) contains a
# small label box ("Rectangle 765") that needs to be narrowed slightly
# and have its first run of text corrected from "SCC" to "SSC".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The whole diagram lives inside one big top-level group shape.
$topGroup = $s.Shapes.Item(1)

# Locate the label box by name inside the (flattened) group item list.
$label = $topGroup.GroupItems.Item("Rectangle 765")

# --- 1. Fix the text: "SCC" -> "SSC" (only the first run) ---------------
$tr = $label.TextFrame.TextRange
$tr.Characters(1, 3).Text = "SSC"

# --- 2. Narrow the shape: <a:ext cx="666" .../> -> cx="660" -------------
# NOTE: Width/Height setters on shapes nested in a group write the raw
# local (child-coordinate) unit straight into the OOXML cx/cy as
# EMU-per-point (1 pt = 12700 EMU), i.e. they do NOT re-apply the
# group's child->parent scale on write (only the getter does). Also,
# the runtime truncates rather than rounds when converting the
# assigned point value back to the integer local unit, so nudge each
# target half a unit up to land squarely on the integer instead of
# one less (e.g. 660 instead of 659).
#
# The text edit above also re-triggers the shape's auto-fit (it has
# <a:spAutoFit/>), which clobbers cy with a stray recalculated value,
# so restore both Width and Height explicitly here to their correct
# local-unit values (660 x 126) after editing the text.
$label.Width = (660.5) / 12700
$label.Height = (126.5) / 12700
